$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("Alternative", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Delete()
    $r.InsertAfter("Changed main")
}
